$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.9720616611150132, 1.1218375607160074, 0.3137193930898504, 2.4882685308903847)"
$ws.Range("C2").Value = "NIG(1.6758215367933418, 1.2554112403714965, 3.3855516397469554, 6.1141186501041265)"
$ws.Range("D2").Value = "NIG(1.088024592827917, 0.6294513452024685, 1.9944303175573586, 3.71509910862388)"
$ws.Range("E2").Value = "EXN(1.7136744594633821, 3.698164852340953, 3.98700188111587)"
